$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
  # Row 15
  $ws.Range("H15").Value = 3090.8096
  $ws.Range("I15").Value = 3090.8096
  $ws.Range("J15").Value = 0
  $ws.Range("K15").Value = 9272.4288
  $ws.Range("L15").Value = 0
  $ws.Range("M15").Value = -9103.4288
  # Row 33
  $ws.Range("H33").Value = 319.92856
  $ws.Range("I33").Value = 267.76923
  $ws.Range("J33").Value = 998
  $ws.Range("K33").Value = 267.76923
  $ws.Range("L33").Value = 998
  $ws.Range("M33").Value = -38.76922999999999
  $ws.Range("N33").Value = -1456
  # Row 43
  $ws.Range("H43").Value = 1500.5
  $ws.Range("I43").Value = 999
  $ws.Range("J43").Value = 2002
  $ws.Range("K43").Value = 999
  $ws.Range("L43").Value = 2002
  $ws.Range("M43").Value = -930
  $ws.Range("N43").Value = -2140
  # Row 97
  $ws.Range("H97").Value = 2151.3333
  $ws.Range("I97").Value = 0
  $ws.Range("J97").Value = 2151.3333
  $ws.Range("K97").Value = 0
  $ws.Range("L97").Value = 6453.999899999999
  $ws.Range("N97").Value = -7445.999899999999
  # Row 111
  $ws.Range("H111").Value = 3029.111
  $ws.Range("I111").Value = 4515
  $ws.Range("J111").Value = 1543.2222
  $ws.Range("K111").Value = 13545
  $ws.Range("L111").Value = 4629.6666
  $ws.Range("M111").Value = -10478
  $ws.Range("N111").Value = -10763.6666
  # Row 113
  $ws.Range("H113").Value = 27780878
  $ws.Range("I113").Value = 33335852
  $ws.Range("J113").Value = 6006
  $ws.Range("K113").Value = 33335852
  $ws.Range("L113").Value = 6006
  $ws.Range("M113").Value = -33332598
  $ws.Range("N113").Value = -12514
  # Row 116
  $ws.Range("H116").Value = 12999.4
  $ws.Range("I116").Value = 2998.5
  $ws.Range("J116").Value = 19666.666
  $ws.Range("K116").Value = 2998.5
  $ws.Range("L116").Value = 19666.666
  $ws.Range("M116").Value = 443.5
  $ws.Range("N116").Value = -26550.666
  # Row 137
  $ws.Range("H137").Value = 3011.5908
  $ws.Range("I137").Value = 3242
  $ws.Range("J137").Value = 2397.1667
  $ws.Range("K137").Value = 9726
  $ws.Range("L137").Value = 7191.500100000001
  $ws.Range("M137").Value = -7176
  $ws.Range("N137").Value = -12291.5001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
  # Row 2
  $ws.Range("H2").Value = 1069.5186
  $ws.Range("I2").Value = 954.55
  $ws.Range("J2").Value = 1398
  $ws.Range("K2").Value = 954.55
  $ws.Range("L2").Value = 1398
  $ws.Range("M2").Value = -841.55
  $ws.Range("N2").Value = -1624
  # Row 35
  $ws.Range("H35").Value = 2000
  $ws.Range("I35").Value = 2000
  $ws.Range("J35").Value = 0
  $ws.Range("K35").Value = 2000
  $ws.Range("L35").Value = 0
  $ws.Range("M35").Value = -1594
  # Row 45
  $ws.Range("H45").Value = 62631.145
  $ws.Range("I45").Value = 72570
  $ws.Range("J45").Value = 2998
  $ws.Range("K45").Value = 72570
  $ws.Range("L45").Value = 2998
  $ws.Range("M45").Value = -72193
  $ws.Range("N45").Value = -3752
  # Row 116
  $ws.Range("H116").Value = 1069.5186
  $ws.Range("I116").Value = 954.55
  $ws.Range("J116").Value = 1398
  $ws.Range("K116").Value = 954.55
  $ws.Range("L116").Value = 1398
  $ws.Range("M116").Value = 1339.45
  $ws.Range("N116").Value = -5986

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
  # Row 3
  $ws.Range("H3").Value = 1069.5186
  $ws.Range("I3").Value = 954.55
  $ws.Range("J3").Value = 1398
  $ws.Range("K3").Value = 954.55
  $ws.Range("L3").Value = 1398
  $ws.Range("M3").Value = -840.55
  $ws.Range("N3").Value = -1626
  # Row 20
  $ws.Range("H20").Value = 31255752
  $ws.Range("I20").Value = 62509130
  $ws.Range("J20").Value = 2375.75
  $ws.Range("K20").Value = 62509130
  $ws.Range("L20").Value = 2375.75
  $ws.Range("M20").Value = -62508883
  $ws.Range("N20").Value = -2869.75
  # Row 80
  $ws.Range("H80").Value = 125000344
  $ws.Range("I80").Value = 250000240
  $ws.Range("J80").Value = 447
  $ws.Range("K80").Value = 250000240
  $ws.Range("L80").Value = 447
  $ws.Range("M80").Value = -249999242
  $ws.Range("N80").Value = -2443
  # Row 83
  $ws.Range("H83").Value = 125000344
  $ws.Range("I83").Value = 250000240
  $ws.Range("J83").Value = 447
  $ws.Range("K83").Value = 1250001200
  $ws.Range("L83").Value = 2235
  $ws.Range("M83").Value = -1249996208
  $ws.Range("N83").Value = -12219
  # Row 94
  $ws.Range("H94").Value = 86961420
  $ws.Range("I94").Value = 100005480
  $ws.Range("J94").Value = 1029.6666
  $ws.Range("K94").Value = 100005480
  $ws.Range("L94").Value = 1029.6666
  $ws.Range("M94").Value = -100005029
  $ws.Range("N94").Value = -1931.6666
  # Row 107
  $ws.Range("H107").Value = 1710676.2
  $ws.Range("I107").Value = 2080211
  $ws.Range("J107").Value = 1577.875
  $ws.Range("K107").Value = 2080211
  $ws.Range("L107").Value = 1577.875
  $ws.Range("M107").Value = -2078291
  $ws.Range("N107").Value = -5417.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
  # Row 28
  $ws.Range("H28").Value = 17000
  $ws.Range("I28").Value = 17000
  $ws.Range("J28").Value = 0
  $ws.Range("K28").Value = 17000
  $ws.Range("L28").Value = 0
  $ws.Range("M28").Value = -16755
  $ws.Range("N28").ClearContents()
  # Row 94
  $ws.Range("H94").Value = 637.5
  $ws.Range("I94").Value = 480.8
  $ws.Range("J94").Value = 749.4286
  $ws.Range("K94").Value = 480.8
  $ws.Range("L94").Value = 749.4286
  $ws.Range("M94").Value = -29.80000000000001
  $ws.Range("N94").Value = -1651.4286
  # Row 132
  $ws.Range("H132").Value = 3904.3333
  $ws.Range("I132").Value = 3189.2917
  $ws.Range("J132").Value = 5811.1113
  $ws.Range("K132").Value = 9567.875100000001
  $ws.Range("L132").Value = 17433.3339
  $ws.Range("M132").Value = -7037.875100000001
  $ws.Range("N132").Value = -22493.3339

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
  # Row 74
  $ws.Range("H74").Value = 21995.6
  $ws.Range("I74").Value = 25006.5
  $ws.Range("J74").Value = 21242.875
  $ws.Range("K74").Value = 75019.5
  $ws.Range("L74").Value = 63728.625
  $ws.Range("M74").Value = -73958.5
  $ws.Range("N74").Value = -65850.625
  # Row 77
  $ws.Range("H77").Value = 21995.6
  $ws.Range("I77").Value = 25006.5
  $ws.Range("J77").Value = 21242.875
  $ws.Range("K77").Value = 225058.5
  $ws.Range("L77").Value = 191185.875
  $ws.Range("M77").Value = -219754.5
  $ws.Range("N77").Value = -201793.875
  # Row 129
  $ws.Range("H129").Value = 112520
  $ws.Range("I129").Value = 0
  $ws.Range("J129").Value = 112520
  $ws.Range("K129").Value = 0
  $ws.Range("L129").Value = 337560
  $ws.Range("N129").Value = -347560
  # Row 140
  $ws.Range("H140").Value = 1270
  $ws.Range("I140").Value = 1270
  $ws.Range("J140").Value = 0
  $ws.Range("K140").Value = 3810
  $ws.Range("L140").Value = 0
  $ws.Range("M140").Value = 1370

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
  # Row 7
  $ws.Range("H7").Value = 15000
  $ws.Range("I7").Value = 0
  $ws.Range("J7").Value = 15000
  $ws.Range("K7").Value = 0
  $ws.Range("L7").Value = 15000
  $ws.Range("N7").Value = -15224
  # Row 8
  $ws.Range("H8").Value = 15000
  $ws.Range("I8").Value = 0
  $ws.Range("J8").Value = 15000
  $ws.Range("K8").Value = 0
  $ws.Range("L8").Value = 15000
  $ws.Range("N8").Value = -15278
  # Row 22
  $ws.Range("H22").Value = 166.6
  $ws.Range("I22").Value = 400
  $ws.Range("J22").Value = 11
  $ws.Range("K22").Value = 400
  $ws.Range("L22").Value = 11
  $ws.Range("M22").Value = 129
  $ws.Range("N22").Value = -1069
  # Row 43
  $ws.Range("H43").Value = 5500
  $ws.Range("I43").Value = 5500
  $ws.Range("J43").Value = 0
  $ws.Range("K43").Value = 5500
  $ws.Range("L43").Value = 0
  $ws.Range("M43").Value = -5349
  $ws.Range("N43").ClearContents()
  # Row 113
  $ws.Range("H113").Value = 2189.45
  $ws.Range("I113").Value = 2164.182
  $ws.Range("J113").Value = 2220.3333
  $ws.Range("K113").Value = 2164.182
  $ws.Range("L113").Value = 2220.3333
  $ws.Range("M113").Value = 5.818000000000211
  $ws.Range("N113").Value = -6560.3333
  # Row 132
  $ws.Range("H132").Value = 2228.4827
  $ws.Range("I132").Value = 1783.5333
  $ws.Range("J132").Value = 2705.2144
  $ws.Range("K132").Value = 5350.5999
  $ws.Range("L132").Value = 8115.6432
  $ws.Range("M132").Value = -2820.5999
  $ws.Range("N132").Value = -13175.6432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
  # Row 7
  $ws.Range("H7").Value = 2628.7273
  $ws.Range("I7").Value = 1903.1666
  $ws.Range("J7").Value = 3499.4
  $ws.Range("K7").Value = 1903.1666
  $ws.Range("L7").Value = 3499.4
  $ws.Range("M7").Value = -1791.1666
  $ws.Range("N7").Value = -3723.4
  # Row 24
  $ws.Range("H24").Value = 0
  $ws.Range("I24").Value = 0
  $ws.Range("J24").Value = 0
  $ws.Range("K24").Value = 0
  $ws.Range("L24").Value = 0
  $ws.Range("M24").ClearContents()
  # Row 40
  $ws.Range("H40").Value = 10888.889
  $ws.Range("I40").Value = 10888.889
  $ws.Range("J40").Value = 0
  $ws.Range("K40").Value = 10888.889
  $ws.Range("L40").Value = 0
  $ws.Range("M40").Value = -10752.889
  # Row 61
  $ws.Range("H61").Value = 1500
  $ws.Range("I61").Value = 1500
  $ws.Range("J61").Value = 0
  $ws.Range("K61").Value = 1500
  $ws.Range("L61").Value = 0
  $ws.Range("M61").Value = -1298
  # Row 113
  $ws.Range("H113").Value = 1500
  $ws.Range("I113").Value = 1500
  $ws.Range("J113").Value = 0
  $ws.Range("K113").Value = 1500
  $ws.Range("L113").Value = 0
  $ws.Range("M113").Value = 670
  # Row 122
  $ws.Range("H122").Value = 6142.091
  $ws.Range("I122").Value = 5326.75
  $ws.Range("J122").Value = 7120.5
  $ws.Range("K122").Value = 15980.25
  $ws.Range("L122").Value = 21361.5
  $ws.Range("M122").Value = -13530.25
  $ws.Range("N122").Value = -26261.5
  # Row 126
  $ws.Range("H126").Value = 2628.7273
  $ws.Range("I126").Value = 1903.1666
  $ws.Range("J126").Value = 3499.4
  $ws.Range("K126").Value = 5709.4998
  $ws.Range("L126").Value = 10498.2
  $ws.Range("M126").Value = -3239.4998
  $ws.Range("N126").Value = -15438.2
  # Row 136
  $ws.Range("H136").Value = 7945.222
  $ws.Range("I136").Value = 8300.4
  $ws.Range("J136").Value = 7501.25
  $ws.Range("K136").Value = 24901.2
  $ws.Range("L136").Value = 22503.75
  $ws.Range("M136").Value = -22351.2
  $ws.Range("N136").Value = -27603.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
  # Row 113
  $ws.Range("H113").Value = 582
  $ws.Range("I113").Value = 582
  $ws.Range("J113").Value = 0
  $ws.Range("K113").Value = 1746
  $ws.Range("L113").Value = 0
  $ws.Range("M113").Value = 424
  # Row 122
  $ws.Range("H122").Value = 16668908
  $ws.Range("I122").Value = 2363.1
  $ws.Range("J122").Value = 50002000
  $ws.Range("K122").Value = 7089.299999999999
  $ws.Range("L122").Value = 150006000
  $ws.Range("M122").Value = -4639.299999999999
  $ws.Range("N122").Value = -150010900
  # Row 136
  $ws.Range("H136").Value = 3322.75
  $ws.Range("I136").Value = 2856.8
  $ws.Range("J136").Value = 4099.3335
  $ws.Range("K136").Value = 8570.400000000001
  $ws.Range("L136").Value = 12298.0005
  $ws.Range("M136").Value = -6020.400000000001
  $ws.Range("N136").Value = -17398.0005
